$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price (D) and Volume(1h) (E) columns are stored as plain text in the workbook
# (inline strings), even though many values look numeric (e.g. "0.9956",
# "28.907.13"). Writing such strings straight into .Value would make Excel silently
# reinterpret them as real numbers, changing their stored representation (precision,
# trailing zeros, thousands separators, etc.). To avoid that we temporarily mark the
# whole D2:E51 data block as Text format ("@") before writing the new values, then
# clear that temporary formatting again afterwards so the cells end up with the same
# (default) style as before - only their text content changes.
$dataRange = $ws.Range("D2:E51")
$dataRange.NumberFormat = "@"

$ws.Range("D2").Value = '28.907.13'
$ws.Range("E2").Value = '  -0.43%  '
$ws.Range("D3").Value = '1.822.55'
$ws.Range("E3").Value = '  -0.53%  '
$ws.Range("D4").Value = '0.9956'
$ws.Range("E4").Value = '  -0.37%  '
$ws.Range("D5").Value = '242.98'
$ws.Range("E5").Value = '  +0.64%  '
$ws.Range("D6").Value = '0.6292'
$ws.Range("E6").Value = '  +0.19%  '
$ws.Range("D7").Value = '0.9953'
$ws.Range("E7").Value = '  -0.52%  '
$ws.Range("D8").Value = '0.07450'
$ws.Range("E8").Value = '  -1.71%  '
$ws.Range("D9").Value = '0.2930'
$ws.Range("E9").Value = '  +0.39%  '
$ws.Range("D10").Value = '22.99'
$ws.Range("E10").Value = '  +0.67%  '
$ws.Range("D11").Value = '0.07673'
$ws.Range("E11").Value = '  +0.33%  '
$ws.Range("D12").Value = '1.822.35'
$ws.Range("E12").Value = '  -0.52%  '
$ws.Range("D13").Value = '4.974'
$ws.Range("E13").Value = '  +0.40%  '
$ws.Range("D14").Value = '0.6658'
$ws.Range("E14").Value = '  -0.09%  '
$ws.Range("D15").Value = '82.79'
$ws.Range("E15").Value = '  +0.34%  '
$ws.Range("D16").Value = '0.000009683'
$ws.Range("E16").Value = '  +1.89%  '
$ws.Range("D17").Value = '6.010'
$ws.Range("E17").Value = '  +0.35%  '
$ws.Range("D18").Value = '28.927.07'
$ws.Range("E18").Value = '  -0.31%  '
$ws.Range("D19").Value = '12.52'
$ws.Range("E19").Value = '  +1.51%  '
$ws.Range("D20").Value = '224.92'
$ws.Range("E20").Value = '  -0.08%  '
$ws.Range("D21").Value = '0.9951'
$ws.Range("E21").Value = '  -0.47%  '
$ws.Range("D22").Value = '7.105'
$ws.Range("E22").Value = '  -1.42%  '
$ws.Range("D23").Value = '0.9977'
$ws.Range("E23").Value = '  -0.29%  '
$ws.Range("D24").Value = '160.00'
$ws.Range("E24").Value = '  -0.11%  '
$ws.Range("E25").Value = '  +3.20%  '
$ws.Range("D26").Value = '8.480'
$ws.Range("E26").Value = '  +0.67%  '
$ws.Range("D27").Value = '17.84'
$ws.Range("E27").Value = '  +0.06%  '
$ws.Range("E28").Value = '  +0.26%  '
$ws.Range("D29").Value = '4.109'
$ws.Range("E29").Value = '  +1.13%  '
$ws.Range("D30").Value = '4.041'
$ws.Range("E30").Value = '  +0.17%  '
$ws.Range("D31").Value = '0.05441'
$ws.Range("E31").Value = '  +4.51%  '
$ws.Range("D32").Value = '1.196'
$ws.Range("E32").Value = '  -0.20%  '
$ws.Range("D33").Value = '1.852'
$ws.Range("E33").Value = '  +0.27%  '
$ws.Range("D34").Value = '0.7415'
$ws.Range("E34").Value = '  +1.28%  '
$ws.Range("D35").Value = '1.133'
$ws.Range("E35").Value = '  -1.53%  '
$ws.Range("D36").Value = '2.602'
$ws.Range("E36").Value = '  +0.43%  '
$ws.Range("D37").Value = '1.237.14'
$ws.Range("E37").Value = '  -2.54%  '
$ws.Range("D38").Value = '2.734'
$ws.Range("E38").Value = '  -0.88%  '
$ws.Range("D39").Value = '0.01773'
$ws.Range("E39").Value = '  -0.91%  '
$ws.Range("D40").Value = '6.670'
$ws.Range("E40").Value = '  +1.96%  '
$ws.Range("D41").Value = '0.8987'
$ws.Range("E41").Value = '  +0.89%  '
$ws.Range("D42").Value = '0.9953'
$ws.Range("E42").Value = '  -0.50%  '
$ws.Range("D43").Value = '101.07'
$ws.Range("D44").Value = '1.971.52'
$ws.Range("E44").Value = '  -0.31%  '
$ws.Range("D47").Value = '0.5060'
$ws.Range("E47").Value = '  -1.05%  '
$ws.Range("D48").Value = '0.4036'
$ws.Range("E48").Value = '  +1.44%  '
$ws.Range("D49").Value = '0.07408'
$ws.Range("E49").Value = '  +4.88%  '
$ws.Range("D50").Value = '8.962'
$ws.Range("E50").Value = '  +1.14%  '
$ws.Range("D51").Value = '1.657'
$ws.Range("E51").Value = '  +1.21%  '

# Row 45 <-> Row 46 content swap: BabyDogeCoin and Aave switch ranking positions,
# each bringing its own updated price/volume data.
$ws.Range("B45").Value = 'Aave'
$ws.Range("C45").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D45").Value = '64.89'
$ws.Range("E45").Value = '  +0.25%  '
$ws.Range("B46").Value = 'BabyDogeCoin'
$ws.Range("C46").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D46").Value = '0.00000000122'
$ws.Range("E46").Value = '  +1.82%  '

# Remove the temporary Text formatting so the cells keep their original (default) style.
$dataRange.ClearFormats()
